# Weekly refresh of the "Fruta, Vega Modelo de Temuco - Níspero" price sheet.
# The reported rows (2-13, 16-20) get reshuffled: each row's Fecha / Calidad /
# Volumen / Precio mínimo / Precio máximo / Precio promedio ponderado /
# Unidad de comercialización / Origen / Precio $/Kg / Kg-unidad now come from
# a different source row in the old sheet (rows 14-15 are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# destination row -> source row (values to copy from, read from the ORIGINAL
# sheet before any writes happen)
$mapping = @{
    2  = 7
    3  = 11
    4  = 2
    5  = 19
    6  = 3
    7  = 18
    8  = 10
    9  = 16
    10 = 17
    11 = 6
    12 = 13
    13 = 20
    16 = 8
    17 = 9
    18 = 4
    19 = 5
    20 = 12
}

# Snapshot every source row's values first so writes to one row never
# clobber data still needed as a source for another row.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $src = $mapping[$row]
    if (-not $snapshot.ContainsKey($src)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$src").Value2
        }
        $snapshot[$src] = $rowVals
    }
}

# Now write the snapshotted values into their destination rows.
foreach ($row in $mapping.Keys) {
    $src = $mapping[$row]
    $rowVals = $snapshot[$src]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
